$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 3464.3333
$ws.Range("I62").Value = 1696.5
$ws.Range("K62").Value = 1696.5
$ws.Range("M62").Value = -1072.5

$ws.Range("H65").Value = 3464.3333
$ws.Range("I65").Value = 1696.5
$ws.Range("K65").Value = 8482.5
$ws.Range("M65").Value = -5362.5

$ws.Range("H111").Value = 1326.65
$ws.Range("I111").Value = 1407.2142
$ws.Range("K111").Value = 4221.642599999999
$ws.Range("M111").Value = -1154.642599999999

$ws.Range("H112").Value = 650327.4399999999
$ws.Range("I112").Value = 566.3333
$ws.Range("J112").Value = 711242.5
$ws.Range("K112").Value = 1698.9999
$ws.Range("L112").Value = 2133727.5
$ws.Range("M112").Value = -590.9999
$ws.Range("N112").Value = -2135943.5

$ws.Range("H138").Value = 2863.7112
$ws.Range("I138").Value = 1328.3667
$ws.Range("J138").Value = 3631.3833
$ws.Range("K138").Value = 3985.1001
$ws.Range("L138").Value = 10894.1499
$ws.Range("M138").Value = 1154.8999
$ws.Range("N138").Value = -21174.1499

$ws.Range("H141").Value = 3531.5454
$ws.Range("I141").Value = 1672.9333
$ws.Range("J141").Value = 7514.2856
$ws.Range("K141").Value = 5018.7999
$ws.Range("L141").Value = 22542.8568
$ws.Range("M141").Value = 161.2001
$ws.Range("N141").Value = -32902.8568

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2954.1785
$ws.Range("I61").Value = 2562.5789
$ws.Range("J61").Value = 3780.889
$ws.Range("K61").Value = 2562.5789
$ws.Range("L61").Value = 3780.889
$ws.Range("M61").Value = -2350.5789
$ws.Range("N61").Value = -4204.889

$ws.Range("H105").Value = 38499.75
$ws.Range("J105").Value = 38499.75
$ws.Range("L105").Value = 38499.75
$ws.Range("N105").Value = -45487.75

$ws.Range("H136").Value = 2954.1785
$ws.Range("I136").Value = 2562.5789
$ws.Range("J136").Value = 3780.889
$ws.Range("K136").Value = 7687.736699999999
$ws.Range("L136").Value = 11342.667
$ws.Range("M136").Value = -5137.736699999999
$ws.Range("N136").Value = -16442.667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H106").Value = 26671
$ws.Range("J106").Value = 26671
$ws.Range("L106").Value = 26671
$ws.Range("N106").Value = -29195

$ws.Range("H134").Value = 3909.6223
$ws.Range("I134").Value = 4162.515
$ws.Range("J134").Value = 3214.1667
$ws.Range("K134").Value = 12487.545
$ws.Range("L134").Value = 9642.500100000001
$ws.Range("M134").Value = -9952.545000000002
$ws.Range("N134").Value = -14712.5001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 25
$ws.Range("I4").Value = 25
$ws.Range("K4").Value = 75
$ws.Range("M4").Value = 37

$ws.Range("H5").Value = 884.37933
$ws.Range("I5").Value = 536.5
$ws.Range("J5").Value = 1453.6364
$ws.Range("K5").Value = 1609.5
$ws.Range("L5").Value = 4360.9092
$ws.Range("M5").Value = -1497.5
$ws.Range("N5").Value = -4584.9092

$ws.Range("H82").Value = 7232.25
$ws.Range("I82").Value = 1000
$ws.Range("J82").Value = 7798.8184
$ws.Range("K82").Value = 3000
$ws.Range("L82").Value = 23396.4552
$ws.Range("M82").Value = -2594
$ws.Range("N82").Value = -24208.4552

$ws.Range("H85").Value = 7232.25
$ws.Range("I85").Value = 1000
$ws.Range("J85").Value = 7798.8184
$ws.Range("K85").Value = 3000
$ws.Range("L85").Value = 23396.4552
$ws.Range("M85").Value = -1596
$ws.Range("N85").Value = -26204.4552

$ws.Range("H92").Value = 567
$ws.Range("I92").Value = 483.875
$ws.Range("J92").Value = 700
$ws.Range("K92").Value = 1451.625
$ws.Range("L92").Value = 2100
$ws.Range("M92").Value = -203.625
$ws.Range("N92").Value = -4596

$ws.Range("H131").Value = 2565026.2
$ws.Range("I131").Value = 16667086
$ws.Range("J131").Value = 1015.4091
$ws.Range("K131").Value = 50001258
$ws.Range("L131").Value = 3046.2273
$ws.Range("M131").Value = -49996218
$ws.Range("N131").Value = -13126.2273

$ws.Range("H134").Value = 1651.1875
$ws.Range("I134").Value = 2236.6667
$ws.Range("J134").Value = 1299.9
$ws.Range("K134").Value = 6710.000100000001
$ws.Range("L134").Value = 3899.7
$ws.Range("M134").Value = -1640.000100000001
$ws.Range("N134").Value = -14039.7

$ws.Range("H135").Value = 884.37933
$ws.Range("I135").Value = 536.5
$ws.Range("J135").Value = 1453.6364
$ws.Range("K135").Value = 4828.5
$ws.Range("L135").Value = 13082.7276
$ws.Range("M135").Value = -2293.5
$ws.Range("N135").Value = -18152.7276

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H105").Value = 30000
$ws.Range("J105").Value = 30000
$ws.Range("L105").Value = 30000
$ws.Range("N105").Value = -36988

$ws.Range("H117").Value = 27000
$ws.Range("J117").Value = 27000
$ws.Range("L117").Value = 27000
$ws.Range("N117").Value = -33884

$ws.Range("H132").Value = 2307.5098
$ws.Range("I132").Value = 1741.3684
$ws.Range("J132").Value = 3962.3845
$ws.Range("K132").Value = 5224.1052
$ws.Range("L132").Value = 11887.1535
$ws.Range("M132").Value = -2694.1052
$ws.Range("N132").Value = -16947.1535

$ws.Range("H135").Value = 27260
$ws.Range("J135").Value = 27260
$ws.Range("L135").Value = 27260
$ws.Range("N135").Value = -37400

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H92").Value = 28000
$ws.Range("J92").Value = 28000
$ws.Range("L92").Value = 28000
$ws.Range("N92").Value = -32992

$ws.Range("H105").Value = 39999
$ws.Range("J105").Value = 39999
$ws.Range("L105").Value = 39999
$ws.Range("N105").Value = -46987

$ws.Range("H122").Value = 3542.9473
$ws.Range("I122").Value = 3663.0344
$ws.Range("J122").Value = 3156
$ws.Range("K122").Value = 10989.1032
$ws.Range("L122").Value = 9468
$ws.Range("M122").Value = -8539.1032
$ws.Range("N122").Value = -14368

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 64812.223
$ws.Range("J46").Value = 64812.223
$ws.Range("L46").Value = 64812.223
$ws.Range("N46").Value = -65274.223

$ws.Range("H81").Value = 32259660
$ws.Range("I81").Value = 58824704
$ws.Range("J81").Value = 2106.8572
$ws.Range("K81").Value = 117649408
$ws.Range("L81").Value = 4213.7144
$ws.Range("M81").Value = -117648347
$ws.Range("N81").Value = -6335.7144

$ws.Range("H82").Value = 23857.143
$ws.Range("J82").Value = 23857.143
$ws.Range("L82").Value = 23857.143
$ws.Range("N82").Value = -24623.143

$ws.Range("H84").Value = 32259660
$ws.Range("I84").Value = 58824704
$ws.Range("J84").Value = 2106.8572
$ws.Range("K84").Value = 588247040
$ws.Range("L84").Value = 21068.572
$ws.Range("M84").Value = -588241736
$ws.Range("N84").Value = -31676.572

$ws.Range("H85").Value = 23857.143
$ws.Range("J85").Value = 23857.143
$ws.Range("L85").Value = 23857.143
$ws.Range("N85").Value = -26509.143

$ws.Range("H97").Value = 28190.666
$ws.Range("J97").Value = 28190.666
$ws.Range("L97").Value = 28190.666
$ws.Range("N97").Value = -30172.666

$ws.Range("H98").Value = 22590
$ws.Range("J98").Value = 22590
$ws.Range("L98").Value = 22590
$ws.Range("N98").Value = -28580

$ws.Range("H104").Value = 28000
$ws.Range("J104").Value = 28000
$ws.Range("L104").Value = 28000
$ws.Range("N104").Value = -34988

$ws.Range("H110").Value = 34296
$ws.Range("J110").Value = 34296
$ws.Range("L110").Value = 34296
$ws.Range("N110").Value = -42476

$ws.Range("H126").Value = 5524.2593
$ws.Range("I126").Value = 6846.1904
$ws.Range("J126").Value = 897.5
$ws.Range("K126").Value = 20538.5712
$ws.Range("L126").Value = 2692.5
$ws.Range("M126").Value = -18068.5712
$ws.Range("N126").Value = -7632.5

$ws.Range("H132").Value = 3124.7817
$ws.Range("I132").Value = 3528.2856
$ws.Range("J132").Value = 2418.65
$ws.Range("K132").Value = 10584.8568
$ws.Range("L132").Value = 7255.950000000001
$ws.Range("M132").Value = -8054.856800000001
$ws.Range("N132").Value = -12315.95

$ws.Range("H134").Value = 64812.223
$ws.Range("J134").Value = 64812.223
$ws.Range("L134").Value = 194436.669
$ws.Range("N134").Value = -199506.669
